$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Get bill" data rows (product1 / ALL / product2 entries)
$ws.Range("A14").Value = "product1"
$ws.Range("B14").Value = 100
$ws.Range("C14").Value = 1

$ws.Range("A15").Value = "product1"
$ws.Range("B15").Value = 50
$ws.Range("C15").Value = "ALL"

$ws.Range("A16").Value = "product2"
$ws.Range("B16").Value = 100
$ws.Range("C16").Value = "ALL"

$ws.Range("A17").Value = "product2"
$ws.Range("B17").Value = 30
$ws.Range("C17").Value = 1

# Update header/footer font style name (Regular -> Normal)
$ps = $ws.PageSetup
$ps.CenterHeader = '&"Times New Roman,Normal"&12&A'
$ps.CenterFooter = '&"Times New Roman,Normal"&12Page &P'

# Move the active selection to D15
$ws.Range("D15").Select()
